$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = @(0.7435897435897436, 0.5188866799204771, 0.6112412177985949)
    3  = @(0.7450980392156863, 0.5288270377733598, 0.6186046511627907)
    4  = @(0.7429094236047575, 0.538104705102717,  0.6241352805534205)
    5  = @(0.7419064748201439, 0.5467196819085487, 0.629530713468142)
    6  = @(0.7522361359570662, 0.5573227302849569, 0.6402740768937952)
    7  = @(0.7548500881834215, 0.5672630881378397, 0.6477487703367386)
    8  = @(0.7572559366754618, 0.5705765407554672, 0.6507936507936508)
    9  = @(0.7555163283318623, 0.5672630881378397, 0.647993943981832)
    10 = @(0.7563989408649603, 0.5679257786613652, 0.6487509462528388)
    11 = @(0.7559681697612732, 0.5666003976143141, 0.6477272727272727)
}

foreach ($row in $values.Keys) {
    $triplet = $values[$row]
    $ws.Range("E$row").Value = $triplet[0]
    $ws.Range("F$row").Value = $triplet[1]
    $ws.Range("G$row").Value = $triplet[2]
}
